# Weekly refresh of "Vega Modelo de Temuco - Lechuga" data:
# two new daily records are inserted above the existing row 589, pushing
# the rest of the historical rows (589:650) down by two rows (to 591:652).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 589; Excel shifts
# everything below (old rows 589:650) down to 591:652 and carries the
# date-column (D) number format down onto the new rows.
$ws.Rows("589:590").Insert()

$newRow589 = @(10, "Vega Modelo de Temuco", "La Araucanía", 44449, 9, 100112033, "Lechuga", "Conconina(o)", "Primera", 220, 7000, 9000, 7864, "`$/caja 10 unidades", "Región Metropolitana", 786, 10, "Hortaliza")
$newRow590 = @(10, "Vega Modelo de Temuco", "La Araucanía", 44449, 9, 100112033, "Lechuga", "Escarola", "Primera", 550, 10000, 12000, 10909, "`$/caja 15 unidades", "Provincia del Elquí", 727, 15, "Hortaliza")

for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(589, $c).Value = $newRow589[$c - 1]
    $ws.Cells.Item(590, $c).Value = $newRow590[$c - 1]
}
